$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.176865220069885
$ws.Range("B1").Value = 2.416972637176514
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.335665464401245
$ws.Range("E1").Value = 1.203413844108582
